$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 13.98369999999999
$ws.Range("C12").Value = -14.72100000000003
$ws.Range("E12").Value = 12.4845
$ws.Range("E14").Value = 13.62100000000001
$ws.Range("E22").Value = 11.8335
